# actualizacion y subida de script
# Updates the "fecha_insercion" column (F) on the "Productos" sheet so that
# every product row reflects the timestamp of the latest run of the
# scraping script (8/24/2025 2:05:21 AM).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

$ws.Range("F2:F21").Value = "8/24/2025 2:05:21 AM"
